# Reduce disturbance in PLBVF data: shrink the deviation of column C
# (voltage set-point) from 1.0 down to one tenth of its original size.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLBVF")
$ws.Activate()

$ws.Range("C3").Value  = 1.001
$ws.Range("C4").Value  = 1.002
$ws.Range("C5").Value  = 1.003
$ws.Range("C6").Value  = 1.002
$ws.Range("C7").Value  = 1.001
$ws.Range("C9").Value  = 0.999
$ws.Range("C10").Value = 0.998
$ws.Range("C11").Value = 0.997
$ws.Range("C12").Value = 0.998
$ws.Range("C13").Value = 0.999

# Scroll down a row and leave the selection on the last data cell, matching
# the view state recorded after the edit.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
